$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Range("A14").Value = 111867661
$ws.Range("B14").Value = 73692
$ws.Range("D14").Value = 'NT'
$ws.Range("E14").Value = 310
$ws.Range("F14").Value = 'Nordlig nållav'
$ws.Range("G14").Value = 'Chaenotheca laevigata'
$ws.Range("H14").Value = 'Nádv.'
$ws.Range("Q14").Value = 703308.4646664646
$ws.Range("R14").Value = 7299302.011735545

# Row 15
$ws.Range("A15").Value = 111866982
$ws.Range("B15").Value = 90682
$ws.Range("D15").Value = 'NT'
$ws.Range("E15").Value = 2059
$ws.Range("F15").Value = 'Skrovlig taggsvamp'
$ws.Range("G15").Value = 'Hydnellum scabrosum'
$ws.Range("H15").Value = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Range("Q15").Value = 703115.6367589685
$ws.Range("R15").Value = 7299535.421832842

# Row 30
$ws.Range("A30").Value = 111867030
$ws.Range("B30").Value = 95538
$ws.Range("D30").Value = 'LC'
$ws.Range("E30").Value = 221941
$ws.Range("F30").Value = 'Plattlummer'
$ws.Range("G30").Value = 'Lycopodium complanatum'
$ws.Range("H30").Value = 'L.'
$ws.Range("Q30").Value = 703118.8428476704
$ws.Range("R30").Value = 7299507.603590234

# Row 31
$ws.Range("A31").Value = 111867456
$ws.Range("B31").Value = 90652
$ws.Range("D31").Value = 'NT'
$ws.Range("E31").Value = 3100
$ws.Range("F31").Value = 'Talltaggsvamp'
$ws.Range("G31").Value = 'Bankera fuligineoalba'
$ws.Range("H31").Value = '(Schmidt : Fr.) Pouzar'
$ws.Range("Q31").Value = 703128.9005519629
$ws.Range("R31").Value = 7299347.87584792

# Row 32
$ws.Range("A32").Value = 111867682
$ws.Range("B32").Value = 90682
$ws.Range("D32").Value = 'NT'
$ws.Range("E32").Value = 2059
$ws.Range("F32").Value = 'Skrovlig taggsvamp'
$ws.Range("G32").Value = 'Hydnellum scabrosum'
$ws.Range("H32").Value = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Range("Q32").Value = 703310.8095286442
$ws.Range("R32").Value = 7299298.053094583

# Column L (empty placeholder cell) moves from row 32 to row 30
$ws.Range("L32").ClearContents()
$ws.Range("L30").NumberFormat = "General"
$ws.Range("L30").ClearFormats()
